# Add a "2022-Q1" holdings sheet (new quarterly snapshot), inserted right
# before the "总计" summary sheet, and add the corresponding aggregate row
# to "总计".

$wb = $excel.ActiveWorkbook

# =============================================================================
# 1. Create the new "2022-Q1" sheet, placed immediately before "总计"
# =============================================================================
# NOTE: a worksheet handle captured BEFORE an Add() call can resolve to a
# different sheet afterwards (the underlying collection shifts), so the new
# sheet is inserted first, and every handle used below is (re-)fetched by
# name fresh, after all Add()/rename operations are done.
$wb.Worksheets.Add($wb.Worksheets.Item("总计")) | Out-Null
$wb.Worksheets.Item("Sheet1").Name = "2022-Q1"

$refSheet = $wb.Worksheets.Item("2021-Q4")   # structural/style template for the new quarter sheet
$newSheet = $wb.Worksheets.Item("2022-Q1")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$data2022 = @(
    ,@("000628", "大成高新技术产业股票A", "57.69", "85.87", "6.51", "3.7556", 8)
    ,@("160910", "大成创新成长混合(LOF)", "16.65", "85.97", "6.92", "1.1522", 3)
    ,@("010929", "大成核心价值甄选混合A", "9.96", "86.12", "7.99", "0.7958", 5)
    ,@("011066", "大成高新技术产业股票C", "10.31", "85.87", "6.51", "0.6712", 8)
    ,@("008271", "大成优势企业混合A", "11.80", "78.45", "5.42", "0.6396", 7)
    ,@("010846", "南方卓越优选3个月持有期混合A", "26.01", "60.36", "2.33", "0.6060", 7)
    ,@("009069", "大成睿鑫股票A", "4.12", "89.40", "8.76", "0.3609", 1)
    ,@("000029", "富国宏观策略灵活配置混合", "5.97", "90.27", "2.04", "0.1218", 9)
    ,@("008272", "大成优势企业混合C", "2.08", "78.45", "5.42", "0.1127", 7)
    ,@("011367", "创金合信群力一年定期开放混合（MOM）A", "3.65", "74.21", "2.39", "0.0872", 4)
    ,@("010847", "南方卓越优选3个月持有期混合C", "3.42", "60.36", "2.33", "0.0797", 7)
    ,@("010930", "大成核心价值甄选混合C", "0.83", "86.12", "7.99", "0.0663", 5)
    ,@("006199", "长盛同锦研究精选混合", "1.73", "82.48", "3.50", "0.0606", 3)
    ,@("008871", "大成睿裕六个月持有期股票A", "0.57", "92.18", "8.25", "0.0470", 4)
    ,@("001892", "长盛新兴成长主题灵活配置混合", "1.32", "82.10", "3.50", "0.0462", 4)
    ,@("009070", "大成睿鑫股票C", "0.42", "89.40", "8.76", "0.0368", 1)
    ,@("002085", "长盛互联网+主题灵活配置混合", "0.84", "83.97", "3.57", "0.0300", 2)
    ,@("004703", "南方兴盛先锋灵活配置混合", "1.09", "53.08", "2.36", "0.0257", 9)
    ,@("011368", "创金合信群力一年定期开放混合（MOM）C", "0.30", "74.21", "2.39", "0.0072", 4)
    ,@("008872", "大成睿裕六个月持有期股票C", "0.04", "92.18", "8.25", "0.0033", 4)
)

$rowCount = $data2022.Count

# Columns B..G on the data rows hold text that looks numeric ("57.69", …) —
# force text formatting BEFORE assigning so it is not coerced into a Number.
$newSheet.Range("B2:G" + ($rowCount + 1)).NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data2022[$i]
    $newSheet.Cells.Item($r, 1).Value = $i          # A: running index (0-based)
    $newSheet.Cells.Item($r, 2).Value = $row[0]     # B: 基金代码
    $newSheet.Cells.Item($r, 3).Value = $row[1]     # C: 基金名称
    $newSheet.Cells.Item($r, 4).Value = $row[2]     # D: 基金规模
    $newSheet.Cells.Item($r, 5).Value = $row[3]     # E: 股票总仓位
    $newSheet.Cells.Item($r, 6).Value = $row[4]     # F: 仓位占比
    $newSheet.Cells.Item($r, 7).Value = $row[5]     # G: 持有市值(亿元)
    $newSheet.Cells.Item($r, 8).Value = $row[6]     # H: 仓位排名 (number)
}

# --- formatting: mirror the look of the other quarter sheets (e.g. 2021-Q4) -
$refSheet.Range("B1:H1").Copy() | Out-Null
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$refSheet.Range("A2").Copy() | Out-Null
$newSheet.Range("A2:A" + ($rowCount + 1)).PasteSpecial(-4122)   # xlPasteFormats

# B2:G<n> picked up the temporary "@" text number-format as their own style;
# restore the plain (unstyled) look used elsewhere by re-pasting the format
# from an already-unstyled reference cell (keeps the cell TYPE as text).
$refSheet.Range("B2:G2").Copy() | Out-Null
$newSheet.Range("B2:G" + ($rowCount + 1)).PasteSpecial(-4122)   # xlPasteFormats

# =============================================================================
# 2. Update "总计": prepend a row for 2022-Q1, shifting the rest down
# =============================================================================
$totalSheet = $wb.Worksheets.Item("总计")

$totalData = @(
    ,@("2022-Q1", 20, 8.71)
    ,@("2021-Q4", 17, 13.66)
    ,@("2021-Q3", 19, 6.74)
    ,@("2021-Q2", 23, 9.119999999999999)
    ,@("2021-Q1", 16, 5.97)
    ,@("2020-Q4", 17, 5.26)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $r = $i + 2
    $row = $totalData[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i     # A: running index (0-based)
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
}

# A2:A7 all share the same "index column" style already present on A2:A6;
# re-apply it to the (now one-row-longer) range for the appended A7.
$totalSheet.Range("A2").Copy() | Out-Null
$totalSheet.Range("A2:A" + ($totalData.Count + 1)).PasteSpecial(-4122) | Out-Null
